$d = $word.ActiveDocument

# 1) Update the Title paragraph's hanging indent (6237 -> 7088 twips = 311.85pt -> 354.4pt)
$titlePara = $d.Paragraphs(1)
$titlePara.Range.ParagraphFormat.LeftIndent = 354.4
$titlePara.Range.ParagraphFormat.FirstLineIndent = -354.4

# 2) Apply italic + size 14 (w:sz/w:szCs 28 half-points) to the "apiVersion" and
#    ' = "v2"' runs, including the complex-script variants (iCs/szCs), which this
#    object model exposes as the "Bi" (bidi) font members.

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-RunFormatting($rangeStart, $rangeEnd, $text) {
    $rng = $d.Range($rangeStart, $rangeEnd)
    $escaped = $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $ns + '><w:body><w:p><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">' + $escaped + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}

# Locate "apiVersion"
$find1 = $d.Content
$find1.Find.Execute("apiVersion", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start1 = $find1.Start
$end1 = $find1.End

Set-RunFormatting $start1 $end1 "apiVersion"

# Locate ' = "v2"' (re-find after the mutation above, text/positions before it are unaffected)
$find2 = $d.Content
$find2.Find.Execute(" = ""v2""", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start2 = $find2.Start
$end2 = $find2.End

Set-RunFormatting $start2 $end2 " = ""v2"""
